# Update the "Price" (D) and "Volume(1h)" (E) columns for the refreshed
# coin-ranking snapshot. Source values are plain text (not numbers/percents),
# so each write uses a leading apostrophe to force text entry, then clears
# the resulting "quote prefix" formatting so the cell style is left untouched
# (matches the original inlineStr cells, which carry no explicit style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'325.58"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'-3.41%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'44.58"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'1.23%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'5.603"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'-3.01%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.08062"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'-3.22%"
$ws.Range("E5").ClearFormats()
$ws.Range("E6").Value = "'-1.82%"
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'4.304"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'-4.93%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'1.902"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'-3.76%"
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'-6.55%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.9452"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'-0.07%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.1163"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'-6.76%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.1862"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'-5.07%"
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.09872"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-0.77%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.04270"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-5.51%"
$ws.Range("E14").ClearFormats()
$ws.Range("D16").Value = "'0.001282"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'-1.56%"
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'0.04210"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'-4.84%"
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'0.005982"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'-1.59%"
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'3.593"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'2.72%"
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'0.3499"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'-0.34%"
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'8.433"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'-4.11%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'0.1371"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'-0.01%"
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.2613"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'-2.95%"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'0.001244"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'-1.43%"
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'0.004467"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'2.47%"
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'-6.37%"
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'0.0003996"
$ws.Range("D27").ClearFormats()
$ws.Range("D39").Value = "'0.02614"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'-6.68%"
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.05430"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'-6.59%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.007716"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'-2.69%"
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.1396"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'-2.39%"
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.007133"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'-20.52%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.002026"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-4.58%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.008573"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'-14.83%"
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.00007137"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'-2.13%"
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'-0.03%"
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.003669"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'15.14%"
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.002274"
$ws.Range("D49").ClearFormats()
$ws.Range("E50").Value = "'-0.03%"
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = "'-0.03%"
$ws.Range("E51").ClearFormats()
